$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "43 x 97" + [char]11 + "  9    7" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "3|    |"
$t.Cell(1,2).Range.Text = "37 x 59" + [char]11 + "  5    9" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "7|    |"
$t.Cell(1,3).Range.Text = "67 x 51" + [char]11 + "  5    1" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "7|    |"
$t.Cell(2,1).Range.Text = "66 x 59" + [char]11 + "  5    9" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "6|    |"
$t.Cell(2,2).Range.Text = "99 x 99" + [char]11 + "  9    9" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "9|    |"
$t.Cell(2,3).Range.Text = "41 x 11" + [char]11 + "  1    1" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "1|    |"
$t.Cell(3,1).Range.Text = "32 x 47" + [char]11 + "  4    7" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "2|    |"
$t.Cell(3,2).Range.Text = "33 x 38" + [char]11 + "  3    8" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "3|    |"
$t.Cell(3,3).Range.Text = "89 x 17" + [char]11 + "  1    7" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "9|    |"
$t.Cell(4,1).Range.Text = "25 x 63" + [char]11 + "  6    3" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "5|    |"
$t.Cell(4,2).Range.Text = "84 x 31" + [char]11 + "  3    1" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "4|    |"
$t.Cell(4,3).Range.Text = "20 x 13" + [char]11 + "  1    3" + [char]11 + "  ----" + [char]11 + "2|    |" + [char]11 + "0|    |"
$t.Cell(5,1).Range.Text = "83 x 82" + [char]11 + "  8    2" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "3|    |"
$t.Cell(5,2).Range.Text = "98 x 71" + [char]11 + "  7    1" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "8|    |"
$t.Cell(5,3).Range.Text = "68 x 93" + [char]11 + "  9    3" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "8|    |"
